# Update gh-pages output data (generated at 456a3b4)
# Sheet 1 = "展览" (Exhibitions), Sheet 3 = "本地生活" (Local life),
# Sheet 4 = "全部类型" (All categories) - a combined view of every sheet.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item(1)
$wsLocal = $wb.Worksheets.Item(3)
$wsAll = $wb.Worksheets.Item(4)

# ---- Sheet 1: 展览 ----
$wsExpo.Range("F2").Value = 1528
$wsExpo.Range("F3").Value = 869
$wsExpo.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202406/OONRvz5A1719312770502.jpeg"
$wsExpo.Range("F4").Value = 452
$wsExpo.Range("F5").Value = 899
$wsExpo.Range("F6").Value = 510
$wsExpo.Range("F7").Value = 7662
$wsExpo.Range("F10").Value = 1938
$wsExpo.Range("F11").Value = 5541
$wsExpo.Range("F14").Value = 7639
$wsExpo.Range("F15").Value = 9031
$wsExpo.Range("F16").Value = 1148
$wsExpo.Range("F17").Value = 903
$wsExpo.Range("F18").Value = 4455
$wsExpo.Range("F19").Value = 675
$wsExpo.Range("F20").Value = 231
$wsExpo.Range("F22").Value = 286
$wsExpo.Range("F25").Value = 117
$wsExpo.Range("F26").Value = 1661
$wsExpo.Range("F27").Value = 721
$wsExpo.Range("F28").Value = 933
$wsExpo.Range("F30").Value = 1874
$wsExpo.Range("F31").Value = 337
$wsExpo.Range("F32").Value = 2299
$wsExpo.Range("F33").Value = 310
$wsExpo.Range("F34").Value = 112
$wsExpo.Range("F35").Value = 1468
$wsExpo.Range("F38").Value = 796
$wsExpo.Range("F39").Value = 510
$wsExpo.Range("F40").Value = 2975
$wsExpo.Range("F41").Value = 4102
$wsExpo.Range("F43").Value = 40
$wsExpo.Range("F44").Value = 421
$wsExpo.Range("F48").Value = 172
$wsExpo.Range("F49").Value = 4087

# ---- Sheet 3: 本地生活 ----
$wsLocal.Range("F2").Value = 5240

# ---- Sheet 4: 全部类型 ----
$wsAll.Range("F3").Value = 1528
$wsAll.Range("F4").Value = 869
$wsAll.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202406/OONRvz5A1719312770502.jpeg"
$wsAll.Range("F5").Value = 452
$wsAll.Range("F6").Value = 899
$wsAll.Range("F7").Value = 510
$wsAll.Range("F11").Value = 5541
$wsAll.Range("F13").Value = 7639
$wsAll.Range("F15").Value = 1148
$wsAll.Range("F16").Value = 903
$wsAll.Range("F17").Value = 675
$wsAll.Range("F18").Value = 231
$wsAll.Range("F20").Value = 286
$wsAll.Range("F24").Value = 117
$wsAll.Range("F25").Value = 1661
$wsAll.Range("F26").Value = 721
$wsAll.Range("F27").Value = 933
$wsAll.Range("F29").Value = 1874
$wsAll.Range("F30").Value = 337
$wsAll.Range("F31").Value = 2299
$wsAll.Range("F39").Value = 510
$wsAll.Range("F40").Value = 4102
$wsAll.Range("F43").Value = 40
$wsAll.Range("F44").Value = 421
$wsAll.Range("F48").Value = 172
$wsAll.Range("F49").Value = 4087
